$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.856.56'
$ws.Range('E2').Value = '  +5.92%  '
$ws.Range('D3').Value = '3.536.55'
$ws.Range('E3').Value = '  +9.52%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '187.06'
$ws.Range('E5').Value = '  +9.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '553.12'
$ws.Range('E6').Value = '  +5.24%  '
$ws.Range('D7').Value = '3.530.69'
$ws.Range('E7').Value = '  +9.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.609'
$ws.Range('E8').Value = '  +2.65%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.630'
$ws.Range('E10').Value = '  +4.36%  '
$ws.Range('E11').Value = '  +14.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.51'
$ws.Range('E12').Value = '  +2.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000269'
$ws.Range('E13').Value = '  +6.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.34'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = '4.108.72'
$ws.Range('E15').Value = '  +9.69%  '
$ws.Range('D16').Value = '3.540.58'
$ws.Range('E16').Value = '  +9.39%  '
$ws.Range('E17').Value = '  +4.52%  '
$ws.Range('D18').Value = '66.874.71'
$ws.Range('E18').Value = '  +6.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.16'
$ws.Range('E19').Value = '  +5.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.94'
$ws.Range('E20').Value = '  +8.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.992'
$ws.Range('E21').Value = '  +2.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '431.92'
$ws.Range('E22').Value = '  +18.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.09'
$ws.Range('E23').Value = '  +8.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.96'
$ws.Range('E24').Value = '  +4.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.09'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.08'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.89'
$ws.Range('E27').Value = '  +9.54%  '
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.15'
$ws.Range('E29').Value = '  +8.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.10'
$ws.Range('E30').Value = '  +11.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.25'
$ws.Range('E31').Value = '  +6.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '647.74'
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.58'
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.72'
$ws.Range('E34').Value = '  +4.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.112'
$ws.Range('E35').Value = '  +5.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.59'
$ws.Range('E36').Value = '  +4.67%  '
$ws.Range('E37').Value = '  +23.61%  '
$ws.Range('D38').Value = '0.0₃0823'
$ws.Range('E38').Value = '  +15.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.48'
$ws.Range('E39').Value = '  +5.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.390'
$ws.Range('E41').Value = '  +3.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.37'
$ws.Range('E42').Value = '  +14.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '3.045.72'
$ws.Range('E44').Value = '  +5.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.68'
$ws.Range('E45').Value = '  +3.09%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.88'
$ws.Range('E46').Value = '  +11.89%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.34'
$ws.Range('E47').Value = '  +9.46%  '
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('E49').Value = '  +6.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.131'
$ws.Range('E50').Value = '  +4.89%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.64'
$ws.Range('E51').Value = '  +11.22%  '
